# "Adjusted firmware for PIC16F1516."
#
# The PWM-backlight sheet computes PIC PWM-timer settings from two user
# inputs: FOSC (B1, oscillator frequency in MHz) and PR2 (B2, the PWM
# period register value). Re-targeting the firmware to a PIC16F1516
# (16 MHz osc, 8-bit/0x3F period) changes just those two inputs; every
# other cell on all three sheets (PWM-backlight, RC-Backlight,
# RC-Contrast) is a formula that recalculates automatically from them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PWM-backlight")

$ws.Range("B1").Value = 16
$ws.Range("B2").Value = 63
